$d = $word.ActiveDocument

# ==========================================================================
# Part 1: Table of contents page numbers shifted down because content was
# removed elsewhere in the document. Each Contents entry is its own
# paragraph ending in a lone run holding just the page number (last run in
# the paragraph, smaller font size than the tab runs before it), so we can
# safely scope a Find/Replace to that paragraph's Range.
# ==========================================================================
$tocFixes = @(
    @{ Idx = 8;  Old = "3";  New = "2" },
    @{ Idx = 9;  Old = "3";  New = "2" },
    @{ Idx = 10; Old = "4";  New = "3" },
    @{ Idx = 11; Old = "4";  New = "3" },
    @{ Idx = 12; Old = "4";  New = "3" },
    @{ Idx = 13; Old = "5";  New = "4" },
    @{ Idx = 14; Old = "7";  New = "4" },
    @{ Idx = 15; Old = "7";  New = "4" },
    @{ Idx = 16; Old = "7";  New = "4" },
    @{ Idx = 17; Old = "7";  New = "4" },
    @{ Idx = 18; Old = "8";  New = "4" },
    @{ Idx = 19; Old = "8";  New = "4" },
    @{ Idx = 20; Old = "8";  New = "5" },
    @{ Idx = 21; Old = "9";  New = "5" },
    @{ Idx = 22; Old = "9";  New = "5" },
    @{ Idx = 23; Old = "9";  New = "5" },
    @{ Idx = 24; Old = "10"; New = "6" },
    @{ Idx = 25; Old = "10"; New = "6" },
    @{ Idx = 26; Old = "10"; New = "6" },
    @{ Idx = 27; Old = "10"; New = "6" },
    @{ Idx = 28; Old = "10"; New = "6" },
    @{ Idx = 29; Old = "11"; New = "6" }
)

foreach ($fix in $tocFixes) {
    $p = $d.Paragraphs.Item($fix.Idx)
    $r = $p.Range
    [void]$r.Find.Execute($fix.Old, $true, $true, $false, $false, $false, $true, 1, $false, $fix.New, 2)
}

# ==========================================================================
# Part 2: " of tests g" + "enerated can be adjusted in the test " (with a
# _GoBack bookmark sitting between the two runs) need to become a single
# run of text " of tests generated can be adjusted in the test ", with no
# bookmark left behind. The neighbouring runs ("number" before, "source "
# after) happen to share identical run formatting with the two runs being
# merged, so a plain text-splice would cascade-merge them in too. Shield
# them with a throwaway formatting toggle while we splice, then remove it.
# ==========================================================================
$r1 = $d.Content.Duplicate
[void]$r1.Find.Execute(" of tests g")

$numberGuard = $d.Range($r1.Start - 6, $r1.Start)

$r2 = $d.Content.Duplicate
[void]$r2.Find.Execute("enerated can be adjusted in the test ")

$sourceGuard = $d.Range($r2.End, $r2.End + 7)

$numberGuard.Bold = 1
$sourceGuard.Bold = 1

# Splice the two runs' text together into one (use a throwaway trailing
# marker character to force a real content change, since the concatenated
# text is otherwise byte-identical to what's already there and a no-op
# text assignment would not actually touch the run structure).
$spliceStart = $r1.Start
$spliceEnd = $r2.End
$full = $d.Range($spliceStart, $spliceEnd)
$full.Text = " of tests generated can be adjusted in the test #"

$markerFix = $d.Content.Duplicate
[void]$markerFix.Find.Execute(" of tests generated can be adjusted in the test #")
$marker = $d.Range($markerFix.End - 1, $markerFix.End)
$marker.Text = " "

# Remove the guards now that the splice is a single run.
$numberGuard2 = $d.Range($r1.Start - 6, $r1.Start)
$numberGuard2.Bold = 0
$sourceGuard2 = $d.Range($spliceStart, $spliceStart)
$sourceGuard2.Collapse(0)
$sourceGuard3 = $d.Range($r1.Start, $r1.Start)
# Re-find "source " right after the merged run to clear its guard.
$mergedRange = $d.Content.Duplicate
[void]$mergedRange.Find.Execute(" of tests generated can be adjusted in the test ")
$sourceGuardFinal = $d.Range($mergedRange.End, $mergedRange.End + 7)
$sourceGuardFinal.Bold = 0

# ==========================================================================
# Part 3: move the (hidden) _GoBack bookmark from the spot above to the
# point in the K=1 sentence where the last edit actually happened, which
# splits that run in two with the bookmark in between.
# ==========================================================================
$kRange = $d.Content.Duplicate
[void]$kRange.Find.Execute("K=1 – allowed, however the game will alway")
$d.Bookmarks.Add("_GoBack", $d.Range($kRange.End, $kRange.End))
